$d = $word.ActiveDocument

# 1. "Add lumber yard research" -> "Add Workshop"
$d.Content.Find.Execute("Add lumber yard research", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Add Workshop", 2)

# 2. Insert a new list paragraph (text: "Add lumber yard research") right after the
#    paragraph that now reads "Add Workshop". Inserting after that paragraph's range
#    inherits its list formatting (ListParagraph style + numPr).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Add Workshop") {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphAfter()

$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $prevText = $p.Previous().Range.Text.TrimEnd([char]13, [char]7)
        if ($prevText -eq "Add Workshop") {
            $newPara = $p
            break
        }
    }
}
$newPara.Range.Text = "Add lumber yard research"

# 3. "Add lumber yard drop off for wood" -> "Add lanterns/torches"
$d.Content.Find.Execute("Add lumber yard drop off for wood", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Add lanterns/torches", 2)
